# Scheduled runner update: refresh market-price derived columns (H:N) on the
# Phantom_Profits leve-crafting sheets. Source data pulled from the Universalis
# price API; only numeric price/profit cells are touched, nothing structural.

$wb = $excel.ActiveWorkbook

# ================= Sheet: ALC =================
$ws = $wb.Worksheets.Item("ALC")

# Row 9
$ws.Cells.Item(9, 8).Value = 833.2857
$ws.Cells.Item(9, 9).Value = 917.63635
$ws.Cells.Item(9, 10).Value = 524
$ws.Cells.Item(9, 11).Value = 917.63635
$ws.Cells.Item(9, 12).Value = 524
$ws.Cells.Item(9, 13).Value = -748.63635
$ws.Cells.Item(9, 14).Value = -862

# Row 32
$ws.Cells.Item(32, 8).Value = 3662.9
$ws.Cells.Item(32, 9).Value = 2274
$ws.Cells.Item(32, 11).Value = 2274
$ws.Cells.Item(32, 13).Value = -1948

# Row 76
$ws.Cells.Item(76, 8).Value = 66668170
$ws.Cells.Item(76, 10).Value = 0
$ws.Cells.Item(76, 12).Value = 0
$ws.Cells.Item(76, 14).ClearContents()

# Row 79
$ws.Cells.Item(79, 8).Value = 66668170
$ws.Cells.Item(79, 10).Value = 0
$ws.Cells.Item(79, 12).Value = 0
$ws.Cells.Item(79, 14).ClearContents()

# Row 100
$ws.Cells.Item(100, 8).Value = 3843
$ws.Cells.Item(100, 9).Value = 2982.182
$ws.Cells.Item(100, 11).Value = 2982.182
$ws.Cells.Item(100, 13).Value = -2441.182

# Row 101
$ws.Cells.Item(101, 8).Value = 829.3333
$ws.Cells.Item(101, 9).Value = 576.5
$ws.Cells.Item(101, 10).Value = 1335
$ws.Cells.Item(101, 11).Value = 1729.5
$ws.Cells.Item(101, 12).Value = 4005
$ws.Cells.Item(101, 13).Value = -107.5
$ws.Cells.Item(101, 14).Value = -7249

# Row 105
$ws.Cells.Item(105, 8).Value = 12000
$ws.Cells.Item(105, 10).Value = 12000
$ws.Cells.Item(105, 12).Value = 12000
$ws.Cells.Item(105, 14).Value = -18988

# Row 132
$ws.Cells.Item(132, 8).Value = 4174.2666
$ws.Cells.Item(132, 9).Value = 4174.2666
$ws.Cells.Item(132, 11).Value = 12522.7998
$ws.Cells.Item(132, 13).Value = -9992.799800000001

# Row 137
$ws.Cells.Item(137, 8).Value = 2869.3333
$ws.Cells.Item(137, 9).Value = 985.55554
$ws.Cells.Item(137, 10).Value = 3999.6
$ws.Cells.Item(137, 11).Value = 2956.66662
$ws.Cells.Item(137, 12).Value = 11998.8
$ws.Cells.Item(137, 13).Value = -406.66662
$ws.Cells.Item(137, 14).Value = -17098.8


# ================= Sheet: ARM =================
$ws = $wb.Worksheets.Item("ARM")

# Row 45
$ws.Cells.Item(45, 8).Value = 3559.2
$ws.Cells.Item(45, 10).Value = 3250
$ws.Cells.Item(45, 12).Value = 3250
$ws.Cells.Item(45, 14).Value = -4004

# Row 97
$ws.Cells.Item(97, 8).Value = 3240
$ws.Cells.Item(97, 9).Value = 1550
$ws.Cells.Item(97, 11).Value = 1550
$ws.Cells.Item(97, 13).Value = -1054

# Row 102
$ws.Cells.Item(102, 8).Value = 951.8
$ws.Cells.Item(102, 9).Value = 951.8
$ws.Cells.Item(102, 11).Value = 951.8
$ws.Cells.Item(102, 13).Value = 670.2


# ================= Sheet: BSM =================
$ws = $wb.Worksheets.Item("BSM")

# Row 86
$ws.Cells.Item(86, 8).Value = 9595.388999999999
$ws.Cells.Item(86, 9).Value = 3859
$ws.Cells.Item(86, 10).Value = 15331.777
$ws.Cells.Item(86, 11).Value = 3859
$ws.Cells.Item(86, 12).Value = 15331.777
$ws.Cells.Item(86, 13).Value = -2736
$ws.Cells.Item(86, 14).Value = -17577.777

# Row 89
$ws.Cells.Item(89, 8).Value = 9595.388999999999
$ws.Cells.Item(89, 9).Value = 3859
$ws.Cells.Item(89, 10).Value = 15331.777
$ws.Cells.Item(89, 11).Value = 19295
$ws.Cells.Item(89, 12).Value = 76658.88499999999
$ws.Cells.Item(89, 13).Value = -13679
$ws.Cells.Item(89, 14).Value = -87890.88499999999

# Row 94
$ws.Cells.Item(94, 8).Value = 1038.5
$ws.Cells.Item(94, 9).Value = 1077.1666
$ws.Cells.Item(94, 11).Value = 1077.1666
$ws.Cells.Item(94, 13).Value = -626.1666

# Row 110
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()


# ================= Sheet: CRP =================
$ws = $wb.Worksheets.Item("CRP")

# Row 6
$ws.Cells.Item(6, 8).Value = 17502000
$ws.Cells.Item(6, 9).Value = 17502000
$ws.Cells.Item(6, 11).Value = 17502000
$ws.Cells.Item(6, 13).Value = -17501887

# Row 20
$ws.Cells.Item(20, 8).Value = 62484.5
$ws.Cells.Item(20, 10).Value = 62484.5
$ws.Cells.Item(20, 12).Value = 62484.5
$ws.Cells.Item(20, 14).Value = -62956.5

# Row 30
$ws.Cells.Item(30, 8).Value = 62484.5
$ws.Cells.Item(30, 10).Value = 62484.5
$ws.Cells.Item(30, 12).Value = 62484.5
$ws.Cells.Item(30, 14).Value = -62666.5

# Row 58
$ws.Cells.Item(58, 8).Value = 2575.5386
$ws.Cells.Item(58, 9).Value = 1564.6666
$ws.Cells.Item(58, 11).Value = 1564.6666
$ws.Cells.Item(58, 13).Value = -1361.6666

# Row 99
$ws.Cells.Item(99, 8).Value = 2392.4167
$ws.Cells.Item(99, 9).Value = 2556.125
$ws.Cells.Item(99, 10).Value = 2065
$ws.Cells.Item(99, 11).Value = 2556.125
$ws.Cells.Item(99, 12).Value = 2065
$ws.Cells.Item(99, 13).Value = -1058.125
$ws.Cells.Item(99, 14).Value = -5061

# Row 103
$ws.Cells.Item(103, 8).Value = 0
$ws.Cells.Item(103, 9).Value = 0
$ws.Cells.Item(103, 11).Value = 0
$ws.Cells.Item(103, 13).ClearContents()

# Row 126
$ws.Cells.Item(126, 8).Value = 2392.4167
$ws.Cells.Item(126, 9).Value = 2556.125
$ws.Cells.Item(126, 10).Value = 2065
$ws.Cells.Item(126, 11).Value = 7668.375
$ws.Cells.Item(126, 12).Value = 6195
$ws.Cells.Item(126, 13).Value = -5198.375
$ws.Cells.Item(126, 14).Value = -11135

# Row 128
$ws.Cells.Item(128, 8).Value = 62484.5
$ws.Cells.Item(128, 10).Value = 62484.5
$ws.Cells.Item(128, 12).Value = 62484.5
$ws.Cells.Item(128, 14).Value = -72444.5

# Row 134
$ws.Cells.Item(134, 8).Value = 4991.5
$ws.Cells.Item(134, 9).Value = 4991.5
$ws.Cells.Item(134, 11).Value = 14974.5
$ws.Cells.Item(134, 13).Value = -12439.5

# Row 136
$ws.Cells.Item(136, 8).Value = 2575.5386
$ws.Cells.Item(136, 9).Value = 1564.6666
$ws.Cells.Item(136, 11).Value = 4693.9998
$ws.Cells.Item(136, 13).Value = -2143.9998


# ================= Sheet: CUL =================
$ws = $wb.Worksheets.Item("CUL")

# Row 99
$ws.Cells.Item(99, 8).Value = 1352.5
$ws.Cells.Item(99, 9).Value = 1352.5
$ws.Cells.Item(99, 11).Value = 4057.5
$ws.Cells.Item(99, 13).Value = -1811.5

# Row 134
$ws.Cells.Item(134, 8).Value = 1665.3334
$ws.Cells.Item(134, 9).Value = 1665.3334
$ws.Cells.Item(134, 11).Value = 4996.0002
$ws.Cells.Item(134, 13).Value = 73.9997999999996


# ================= Sheet: GSM =================
$ws = $wb.Worksheets.Item("GSM")

# Row 36
$ws.Cells.Item(36, 8).Value = 74130.5
$ws.Cells.Item(36, 10).Value = 11874.5
$ws.Cells.Item(36, 12).Value = 11874.5
$ws.Cells.Item(36, 14).Value = -12844.5

# Row 80
$ws.Cells.Item(80, 8).Value = 3305.8
$ws.Cells.Item(80, 10).Value = 3532
$ws.Cells.Item(80, 12).Value = 3532
$ws.Cells.Item(80, 14).Value = -5528

# Row 83
$ws.Cells.Item(83, 8).Value = 3305.8
$ws.Cells.Item(83, 10).Value = 3532
$ws.Cells.Item(83, 12).Value = 17660
$ws.Cells.Item(83, 14).Value = -27644

# Row 113
$ws.Cells.Item(113, 8).Value = 1037.375
$ws.Cells.Item(113, 9).Value = 971.2857
$ws.Cells.Item(113, 11).Value = 971.2857
$ws.Cells.Item(113, 13).Value = 1198.7143


# ================= Sheet: LTW =================
$ws = $wb.Worksheets.Item("LTW")

# Row 7
$ws.Cells.Item(7, 8).Value = 4999.353
$ws.Cells.Item(7, 9).Value = 4861.8125
$ws.Cells.Item(7, 11).Value = 4861.8125
$ws.Cells.Item(7, 13).Value = -4749.8125

# Row 19
$ws.Cells.Item(19, 8).Value = 3499
$ws.Cells.Item(19, 10).Value = 3499
$ws.Cells.Item(19, 12).Value = 3499
$ws.Cells.Item(19, 14).Value = -3839

# Row 22
$ws.Cells.Item(22, 8).Value = 994.5
$ws.Cells.Item(22, 10).Value = 0
$ws.Cells.Item(22, 12).Value = 0
$ws.Cells.Item(22, 14).ClearContents()

# Row 27
$ws.Cells.Item(27, 8).Value = 994.5
$ws.Cells.Item(27, 10).Value = 0
$ws.Cells.Item(27, 12).Value = 0
$ws.Cells.Item(27, 14).ClearContents()

# Row 40
$ws.Cells.Item(40, 8).Value = 1816
$ws.Cells.Item(40, 9).Value = 1876.75
$ws.Cells.Item(40, 11).Value = 1876.75
$ws.Cells.Item(40, 13).Value = -1740.75

# Row 61
$ws.Cells.Item(61, 8).Value = 1938.4615
$ws.Cells.Item(61, 9).Value = 1890.2
$ws.Cells.Item(61, 11).Value = 1890.2
$ws.Cells.Item(61, 13).Value = -1688.2

# Row 113
$ws.Cells.Item(113, 8).Value = 1938.4615
$ws.Cells.Item(113, 9).Value = 1890.2
$ws.Cells.Item(113, 11).Value = 1890.2
$ws.Cells.Item(113, 13).Value = 279.8

# Row 122
$ws.Cells.Item(122, 8).Value = 3405.1333
$ws.Cells.Item(122, 9).Value = 3320.0908
$ws.Cells.Item(122, 10).Value = 3639
$ws.Cells.Item(122, 11).Value = 9960.2724
$ws.Cells.Item(122, 12).Value = 10917
$ws.Cells.Item(122, 13).Value = -7510.2724
$ws.Cells.Item(122, 14).Value = -15817

# Row 126
$ws.Cells.Item(126, 8).Value = 4999.353
$ws.Cells.Item(126, 9).Value = 4861.8125
$ws.Cells.Item(126, 11).Value = 14585.4375
$ws.Cells.Item(126, 13).Value = -12115.4375

# Row 136
$ws.Cells.Item(136, 8).Value = 3407
$ws.Cells.Item(136, 9).Value = 2195.84
$ws.Cells.Item(136, 11).Value = 6587.52
$ws.Cells.Item(136, 13).Value = -4037.52


# ================= Sheet: WVR =================
$ws = $wb.Worksheets.Item("WVR")

# Row 107
$ws.Cells.Item(107, 8).Value = 1056.3334
$ws.Cells.Item(107, 9).Value = 1168.2858
$ws.Cells.Item(107, 10).Value = 664.5
$ws.Cells.Item(107, 11).Value = 3504.8574
$ws.Cells.Item(107, 12).Value = 1993.5
$ws.Cells.Item(107, 13).Value = -1584.8574
$ws.Cells.Item(107, 14).Value = -5833.5

# Row 110
$ws.Cells.Item(110, 8).Value = 0
$ws.Cells.Item(110, 10).Value = 0
$ws.Cells.Item(110, 12).Value = 0
$ws.Cells.Item(110, 14).ClearContents()

# Row 113
$ws.Cells.Item(113, 8).Value = 558.75
$ws.Cells.Item(113, 9).Value = 440.91666
$ws.Cells.Item(113, 11).Value = 1322.74998
$ws.Cells.Item(113, 13).Value = 847.2500199999999

# Row 126
$ws.Cells.Item(126, 8).Value = 2563.5
$ws.Cells.Item(126, 9).Value = 1377
$ws.Cells.Item(126, 10).Value = 3750
$ws.Cells.Item(126, 11).Value = 4131
$ws.Cells.Item(126, 12).Value = 11250
$ws.Cells.Item(126, 13).Value = -1661
$ws.Cells.Item(126, 14).Value = -16190

# Row 130
$ws.Cells.Item(130, 8).Value = 48990
$ws.Cells.Item(130, 10).Value = 48990
$ws.Cells.Item(130, 12).Value = 48990
$ws.Cells.Item(130, 14).Value = -59030

